$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Rakib): G,H,I were blank, become explicit 0
$ws.Range("G3:I3").Value = 0

# Row 4 (Mahfuz): G,H,I were blank, become explicit 0
$ws.Range("G4:I4").Value = 0

# Row 5 (Himel): G,H,I were blank, become explicit 0
$ws.Range("G5:I5").Value = 0

# Row 6 (Minhaz): G,H become 2.5 (I stays blank)
$ws.Range("G6:H6").Value = 2.5

# Row 7 (Taher): G,H become 2.5 (I stays blank)
$ws.Range("G7:H7").Value = 2.5

# Row 8 (Forhad): G,H become 2.5 (I stays blank)
$ws.Range("G8:H8").Value = 2.5

# Row 9 (Nayem): G,H become 2.5 (I stays blank)
$ws.Range("G9:H9").Value = 2.5

# Row 26 (Minhaz deposit row): G26 0 -> 45 (Minhaz +45)
$ws.Range("G26").Value = 45

# Row 42 (Name header row for bazar log): G42 gets label "Minhaz"
$ws.Range("G42").Value = "Minhaz"

# Row 43 (Bazar TK row): G43 blank -> 545 (Bazar -545, i.e. additional bazar cost)
$ws.Range("G43").Value = 545

# Update the view state to match the authored state (scroll so row 3 is
# the top visible row, and select I3:I5 with I3 as the active cell)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I3:I5").Select()
